$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Test Cases")
$ws2 = $wb.Worksheets.Item("TestData")

# --- Rename "Test Cases" -> "TestCases" ---
$ws1.Name = "TestCases"

# ================= TestData sheet (sheet2) =================

# Row 4: (empty) | java | true | (empty)   -- build before touching C3
$ws2.Range("A2:C2").Copy()
$ws2.Range("A4:C4").PasteSpecial(-4122)
$ws2.Range("B4").Value = "java"
$ws2.Range("C2").Copy()
$ws2.Range("C4").PasteSpecial(-4163)

# Row 5: (empty) | !@[]\#^ (hyperlink) | false | (empty)  -- uses the ORIGINAL false (C3) as source
$ws2.Range("A3:C3").Copy()
$ws2.Range("A5:C5").PasteSpecial(-4122)
$ws2.Range("B5").Value = "!@[]\#^"
$ws2.Range("C3").Copy()
$ws2.Range("C5").PasteSpecial(-4163)

# Now it is safe to change C3 false -> true
$ws2.Range("C2").Copy()
$ws2.Range("C3").PasteSpecial(-4163)

# D2:D5 should have the data-row style (empty cells)
$ws2.Range("A2").Copy()
$ws2.Range("D2:D5").PasteSpecial(-4122)

# Header cells: C1 "Expected result", D1 "Result" -- added last so "Expected result" is the last new shared string
$ws2.Range("C1").Value = "Expected result"
$ws2.Range("D1").Value = "Result"
$ws2.Range("A1:B1").Copy()
$ws2.Range("C1:D1").PasteSpecial(-4122)

# Hyperlink on B5
$ws2.Hyperlinks.Add($ws2.Range("B5"), "http://example.com")

# Column widths (approximate, engine quantizes to 1/6 steps)
$ws2.Columns.Item(1).ColumnWidth = 13.877604166666666
$ws2.Columns.Item(2).ColumnWidth = 9.307291666666666
$ws2.Columns.Item(3).ColumnWidth = 13.736979166666666
$ws2.Columns.Item(4).ColumnWidth = 13.451822916666666

# ================= TestCases sheet (sheet1) =================
# Clear "Passed" result column values (E2:E10)
$ws1.Range("E2:E10").ClearContents()

# ================= Selections / active sheet =================
$ws1.Range("E19").Select()
$ws2.Range("D2").Select()
$ws2.Activate()
